# Add a column for water added by the FlowModel to the Quick Check report,
# and add a new "Baseline_2010_C88+ 10/31/20" data row (row 13) ahead of the
# existing "2010-18" block, pushing the rest of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell (S1) for the added column -----------------------------
# Must be written before any other *new* shared string so it lands at
# sharedStrings index 34 (matching the target workbook).
$ws.Range("S1").Value() = " added reach water (fraction)"
# Match the style used by the rest of row 1's header cells (wrapped text).
$ws.Range("S1").WrapText() = $true

# --- Insert a new row 13, shifting the old rows 13-18 down to 14-19 --------
$ws.Rows.Item(13).Insert()
# The freshly inserted row inherits formatting from the row above (row 12);
# the target workbook has this row completely unstyled, so strip it.
$ws.Range("A13:S13").ClearFormats()

# --- Populate the new row 13 with the Baseline_2010_C88+ data --------------
$ws.Range("A13").Value() = "CW3M"
$ws.Range("B13").Value() = "Baseline_2010_C88+ 10/31/20"
$ws.Range("C13").Value() = 2010
$ws.Range("D13").Value() = 1090.199341
$ws.Range("E13").Value() = 1990.4676509999999
$ws.Range("F13").Value() = 1.255063
$ws.Range("G13").Value() = 327.58108499999997
$ws.Range("H13").Value() = 10.610913999999999
$ws.Range("I13").Value() = 8.8404570000000007
$ws.Range("J13").Value() = 814.49517800000001
$ws.Range("K13").Value() = 93.229797000000005
$ws.Range("L13").Value() = 1305.1243899999999
$ws.Range("M13").Value() = 1201.781982
$ws.Range("N13").Value() = 7126.6015630000002
$ws.Range("O13").Value() = 29450.638672000001
$ws.Range("P13").Value() = 3.3577499999999998
$ws.Range("Q13").Value() = 0.00098200000000000002
$ws.Range("R13").Value() = 2010
$ws.Range("S13").Value() = 0.0011180000000000001

# T13 = Q13 - S13, formatted like the other "mass balance discrepancy"
# fraction cells (0.000000).
$ws.Range("T13").Formula() = "=Q13-S13"
$ws.Range("T13").NumberFormat() = $ws.Range("Q17").NumberFormat()

# --- Column width for the new column T (20) ---------------------------------
$ws.Columns.Item(20).ColumnWidth() = 9.21875

# --- Sheet view tweaks -------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow() = 5
$ws.Range("T13").Select()
